# Improve local DS diagram (slide 4: "Local or Remote Decision Service" / In-Process).
#
# Point values below were chosen so that PowerPoint's Single-precision
# Left/Top/Width/Height properties round-trip to the exact target EMU
# offsets/extents (EMU = floor(float32(points) * 12700)).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Shape 6: "Arrow: Left 5" (id=6) ---------------------------------------
$arrow1 = $s.Shapes.Item(6)
$arrow1.Left   = 358.9179
$arrow1.Top    = 356.9226
$arrow1.Width  = 239.36772
$arrow1.Height = 14.4987

# --- Shape 7: "Arrow: Left 6" (id=7, rotated 180) --------------------------
$arrow2 = $s.Shapes.Item(7)
$arrow2.Left   = 358.918
$arrow2.Top    = 301.89324
$arrow2.Width  = 239.3677
$arrow2.Height = 14.4986

# --- Shape 8: "TextBox 13" (id=14) -> reposition + retext -------------------
$tb13 = $s.Shapes.Item(8)
$tb13.Left   = 420.4874
$tb13.Top    = 268.07
$tb13.Width  = 116.2292
$tb13.Height = 29.0813
$tb13.TextFrame.TextRange.Text = "JavaScript call"

# --- Shape 9: "TextBox 14" (id=15) -> reposition + shrink font -------------
$tb14 = $s.Shapes.Item(9)

# Burn through the slide's shape-id/name counter with throwaway
# add+delete cycles so that the shape we add further down receives the
# same id/name PowerPoint originally assigned it (id=25, "TextBox 24").
for ($i = 1; $i -le 11; $i++) {
    $burn = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $burn.Delete()
}

# Duplicate TextBox 14 *before* moving it, so the new shape inherits all
# of its OOXML boilerplate (noFill, wrap="none", spAutoFit, lstStyle,
# rtlCol="0", dirty="0", ...), then reposition/relabel the duplicate.
$dupRange = $tb14.Duplicate()
$newBox = $dupRange.Item(1)
$newBox.ZOrder(0)   # msoBringToFront -> move to the end of the shape tree
$newBox.Name = "TextBox 24"
$newBox.Left   = 450.59914
$newBox.Top    = 314.1921
$newBox.Width  = 43.9501
$newBox.Height = 24.23442
$newBox.TextFrame.TextRange.Text = "JSON"
$newBox.TextFrame.TextRange.Font.Size = 14

# Now move/shrink the original TextBox 14 into its new spot.
$tb14.Left   = 450.7587
$tb14.Top    = 373.3927
$tb14.Width  = 43.9501
$tb14.Height = 24.23442
$tb14.TextFrame.TextRange.Font.Size = 14
